# Fix "Mortendad" -> "Mortandad" typo, and split the combined
# "Los Alamos and Pajarito Canyons" section into two separate sections.

$wb = $excel.ActiveWorkbook

# --- Sheet "Alluvial for Mapping": fix Watershed column typo ---------------
$wsMap = $wb.Worksheets.Item("Alluvial for Mapping")
$wsMap.Range("P16").Value = "Mortandad"
$wsMap.Range("P17").Value = "Mortandad"
$wsMap.Range("P18").Value = "Mortandad"

# --- Sheet "Alluvial Exhibit": fix header + split canyon section -----------
$wsEx = $wb.Worksheets.Item("Alluvial Exhibit")

# Fix "Mortendad Canyon" -> "Mortandad Canyon" header (row 18)
$wsEx.Range("A18").Value = "Mortandad Canyon"

# Column C needs to be a bit wider for the new "Pajarito Canyon" label.
# (ColumnWidth is expressed in character-width COM units, which Excel then
# rounds to pixel-quantised storage units for the raw <col width=".."/>
# attribute; 12.14 is the COM value that round-trips to a stored width of
# exactly 13, matching the target column width.)
$wsEx.Columns.Item(3).ColumnWidth = 12.14

# "Los Alamos and Pajarito Canyons" (row 22) becomes just "Los Alamos Canyon"
$wsEx.Range("A22").Value = "Los Alamos Canyon"

# Insert a new section header row before the old row 26 (18-MW-18) so the
# Pajarito Canyon wells (18-MW-18, PCAO-8) get their own header, matching
# the style used for the other section headers (e.g. row 22).
$wsEx.Rows.Item(26).Insert()
$wsEx.Range("B26:G26").Clear()
$wsEx.Range("A26:H26").Merge()

$wsEx.Range("A22").Copy()
$wsEx.Range("A26").PasteSpecial(-4122)

$wsEx.Range("H22").Copy()
$wsEx.Range("H26").PasteSpecial(-4122)

$wsEx.Range("A26").Value = "Pajarito Canyon"

$wsEx.Range("B26:G26").ClearFormats()
$wsEx.Range("B26:G26").ClearContents()
$wsEx.Range("B26:G26").Clear()
